$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 126 (shifts existing rows 126:146 down to 127:147)
$ws.Rows.Item(126).Insert()

# Populate the new row 126 with the new weekly price record
$ws.Range("A126").Value = 10
$ws.Range("B126").Value = "Vega Modelo de Temuco"
$ws.Range("C126").Value = "La Araucanía"
$ws.Range("D126").Value = 44776
$ws.Range("D126").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E126").Value = 9
$ws.Range("F126").Value = "Fruta"
$ws.Range("G126").Value = 100104
$ws.Range("H126").Value = "Frutos de pepita"
$ws.Range("I126").Value = 100104001
$ws.Range("J126").Value = "Granada"
$ws.Range("K126").Value = "Wonderfull"
$ws.Range("L126").Value = "Primera"
$ws.Range("M126").Value = 100
$ws.Range("N126").Value = 14000
$ws.Range("O126").Value = 14000
$ws.Range("P126").Value = 14000
$ws.Range("Q126").Value = "$/bandeja 10 kilos granel"
$ws.Range("R126").Value = "Provincia de Limarí"
$ws.Range("S126").Value = 1400
$ws.Range("T126").Value = 10
